$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B:AB content between row 114 and row 115
$row114 = $ws.Range("B114:AB114").Value2
$row115 = $ws.Range("B115:AB115").Value2
$ws.Range("B114:AB114").Value2 = $row115
$ws.Range("B115:AB115").Value2 = $row114

# Swap the B:AB content between row 135 and row 136
$row135 = $ws.Range("B135:AB135").Value2
$row136 = $ws.Range("B136:AB136").Value2
$ws.Range("B135:AB135").Value2 = $row136
$ws.Range("B136:AB136").Value2 = $row135
